$wb = $excel.ActiveWorkbook

$wsVendors  = $wb.Worksheets.Item("Vendors")
$wsClients  = $wb.Worksheets.Item("Clients")
$wsExpenses = $wb.Worksheets.Item("Expenses")

# --- Clients sheet (sheet2): new rows 19-21 -------------------------------
# Row 19 repeats an existing client (James Embrey / his Revature email)
$wsClients.Range("A19").Value = "James Embrey"
$wsClients.Range("B19").Value = "james.embrey@revature.net"

# Row 20 - brand new client "Ai Krasner"
$wsClients.Range("A20").Value = "Ai Krasner"
$wsClients.Range("B20").Value = "ai.krasner@gmail.com"

# --- Expenses sheet (sheet3): new rows 2-3 --------------------------------
$wsExpenses.Range("A2").Value = "James Embrey"
$wsExpenses.Range("B2").Value = "Eggs"
$wsExpenses.Range("C2").Value = 12
$wsExpenses.Range("D2").Value = 19.2

$wsExpenses.Range("A3").Value = "Ai Krasner"
$wsExpenses.Range("B3").Value = "Zelda BOTW (5% Discount)"
$wsExpenses.Range("C3").Value = 2
$wsExpenses.Range("D3").Value = 75.981

# --- back to Clients: row 21 (new email for an existing client "Jake") ----
$wsClients.Range("A21").Value = "Jake"
$wsClients.Range("B21").Value = "jake.peralta@gmail.com"

# --- Vendors sheet (sheet1): restock a couple of existing products -------
$wsVendors.Range("D3").Value = 100
$wsVendors.Range("D6").Value = 100

# --- Vendors sheet (sheet1): new vendors / products, rows 16-21 ----------
$wsVendors.Range("A16").Value = "Soda City"
$wsVendors.Range("B16").Value = "Surge"
$wsVendors.Range("C16").Value = 5.99
$wsVendors.Range("D16").Value = 50

$wsVendors.Range("B17").Value = "Mountain Dew"
$wsVendors.Range("C17").Value = 5.99
$wsVendors.Range("D17").Value = 50

$wsVendors.Range("B18").Value = "Pepsi"
$wsVendors.Range("C18").Value = 7.99
$wsVendors.Range("D18").Value = 50

$wsVendors.Range("A19").Value = "Target"
$wsVendors.Range("B19").Value = "Bicycle"
$wsVendors.Range("C19").Value = 99.99
$wsVendors.Range("D19").Value = 50

$wsVendors.Range("B20").Value = "Washing Machine"
$wsVendors.Range("C20").Value = 149.99
$wsVendors.Range("D20").Value = 50

$wsVendors.Range("B21").Value = "Basket"
$wsVendors.Range("C21").Value = 3.99
$wsVendors.Range("D21").Value = 25

# --- view state: Vendors becomes the active sheet/tab, with D2:D10 selected
# (Expenses' own selection is left untouched - it is not the active sheet any more)
$wsClients.Range("A2:B18").Select()
$wsVendors.Range("D2:D10").Select()
